$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.193.78"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "1.699.81"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.89"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.98"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.265"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0629"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.942.95"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "1.716.42"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.557"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.05"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "250.43"
$ws.Range("E17").Value = "  +5.99%  "
$ws.Range("D18").Value = "28.157.77"
$ws.Range("E18").Value = "  +3.59%  "
$ws.Range("D19").Value = "0.0₃0744"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.54"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.56"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.55"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.33"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.45"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "1.466.61"
$ws.Range("E33").Value = "  -4.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.956"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.593"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0173"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.12"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.61"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "1.849.05"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.797"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +7.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.55"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.94"
$ws.Range("E51").Value = "  -4.27%  "
